# "Add files via upload" -- adds new backlog items (pulled from the Stories
# sheet) to the bottom of the Backlog sheet, and leaves the workbook
# positioned on the Sprint1 tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Backlog: append rows 21-32 (Story ID / Story Name pairs copied from the
# Stories sheet backlog candidates).
# ---------------------------------------------------------------------
$backlog = $wb.Worksheets.Item("Backlog")

$newItems = @(
    @("US15", "Fewer than 15 siblings"),
    @("US12", "Parents not too old"),
    @("US16", "Male last names"),
    @("US18", "Siblings should not marry"),
    @("US21", "Correct gender for role"),
    @("US25", "Unique first names in families"),
    @("US28", "Order siblings by age"),
    @("US31", "List living single"),
    @("US35", "List recent births"),
    @("US36", "List recent deaths"),
    @("US38", "List upcoming birthdays"),
    @("US39", "List upcoming anniversaries")
)

$row = 21
foreach ($item in $newItems) {
    $backlog.Cells.Item($row, 2).Value = $item[0]
    $backlog.Cells.Item($row, 3).Value = $item[1]
    $row = $row + 1
}

# Two of the pasted-in rows (US15 / US16) carry the Stories sheet's wrapped
# description formatting on their trailing (empty) D cell -- mirror that by
# copying the format from the matching Stories rows.
$stories = $wb.Worksheets.Item("Stories")
$stories.Range("C16").Copy()
$backlog.Range("D21").PasteSpecial(-4122)
$stories.Range("C17").Copy()
$backlog.Range("D23").PasteSpecial(-4122)

$backlog.Range("B33").Select()

# ---------------------------------------------------------------------
# Move the active tab to Sprint1 (was Backlog).
# ---------------------------------------------------------------------
$sprint1 = $wb.Worksheets.Item("Sprint1")
$sprint1.Activate()
$sprint1.Range("B24").Select()
